# Auto-generated Excel COM-interop script applying numeric-value updates
# to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR),
# reflecting a refreshed market-price pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3781.7273
$ws.Range("I70").Value = 3574.5
$ws.Range("J70").Value = 3827.7778
$ws.Range("K70").Value = 10723.5
$ws.Range("L70").Value = 11483.3334
$ws.Range("M70").Value = -10453.5
$ws.Range("N70").Value = -12023.3334

$ws.Range("H73").Value = 3781.7273
$ws.Range("I73").Value = 3574.5
$ws.Range("J73").Value = 3827.7778
$ws.Range("K73").Value = 10723.5
$ws.Range("L73").Value = 11483.3334
$ws.Range("M73").Value = -9787.5
$ws.Range("N73").Value = -13355.3334

$ws.Range("H80").Value = 1334.6296
$ws.Range("I80").Value = 1704
$ws.Range("J80").Value = 1117.3529
$ws.Range("K80").Value = 5112
$ws.Range("L80").Value = 3352.0587
$ws.Range("M80").Value = -4114
$ws.Range("N80").Value = -5348.0587

$ws.Range("H83").Value = 1334.6296
$ws.Range("I83").Value = 1704
$ws.Range("J83").Value = 1117.3529
$ws.Range("K83").Value = 15336
$ws.Range("L83").Value = 10056.1761
$ws.Range("M83").Value = -10344
$ws.Range("N83").Value = -20040.1761

$ws.Range("H106").Value = 11291.27
$ws.Range("I106").Value = 4117.3
$ws.Range("K106").Value = 4117.3
$ws.Range("M106").Value = -3486.3

$ws.Range("H132").Value = 11153.792
$ws.Range("I132").Value = 11421.348
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 34264.044
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -31734.044
$ws.Range("N132").Value = -20060

$ws.Range("H137").Value = 33345384
$ws.Range("I137").Value = 55557860
$ws.Range("J137").Value = 26674.5
$ws.Range("K137").Value = 166673580
$ws.Range("L137").Value = 80023.5
$ws.Range("M137").Value = -166671030
$ws.Range("N137").Value = -85123.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 643567.6
$ws.Range("I32").Value = 713975.9
$ws.Range("K32").Value = 713975.9
$ws.Range("M32").Value = -713688.9

$ws.Range("H97").Value = 30304002
$ws.Range("I97").Value = 703.25
$ws.Range("J97").Value = 200002480
$ws.Range("K97").Value = 703.25
$ws.Range("L97").Value = 200002480
$ws.Range("M97").Value = -207.25
$ws.Range("N97").Value = -200003472

$ws.Range("H132").Value = 404351.12
$ws.Range("I132").Value = 439576.7
$ws.Range("J132").Value = 2779.6
$ws.Range("K132").Value = 1318730.1
$ws.Range("L132").Value = 8338.799999999999
$ws.Range("M132").Value = -1316200.1
$ws.Range("N132").Value = -13398.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2413.5588
$ws.Range("I86").Value = 2446.65
$ws.Range("J86").Value = 2366.2856
$ws.Range("K86").Value = 2446.65
$ws.Range("L86").Value = 2366.2856
$ws.Range("M86").Value = -1323.65
$ws.Range("N86").Value = -4612.2856

$ws.Range("H89").Value = 2413.5588
$ws.Range("I89").Value = 2446.65
$ws.Range("J89").Value = 2366.2856
$ws.Range("K89").Value = 12233.25
$ws.Range("L89").Value = 11831.428
$ws.Range("M89").Value = -6617.25
$ws.Range("N89").Value = -23063.428

$ws.Range("H94").Value = 1762.6666
$ws.Range("I94").Value = 1902.174
$ws.Range("J94").Value = 960.5
$ws.Range("K94").Value = 1902.174
$ws.Range("L94").Value = 960.5
$ws.Range("M94").Value = -1451.174
$ws.Range("N94").Value = -1862.5

$ws.Range("H99").Value = 6877.619
$ws.Range("I99").Value = 7523.8887
$ws.Range("K99").Value = 7523.8887
$ws.Range("M99").Value = -6025.8887

$ws.Range("H107").Value = 1738.129
$ws.Range("I107").Value = 1696.0667
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 1696.0667
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 223.9332999999999
$ws.Range("N107").Value = -6840

$ws.Range("H135").Value = 98000
$ws.Range("J135").Value = 98000
$ws.Range("L135").Value = 98000
$ws.Range("N135").Value = -108140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3474562.5
$ws.Range("I31").Value = 5052826.5
$ws.Range("J31").Value = 2382.2
$ws.Range("K31").Value = 5052826.5
$ws.Range("L31").Value = 2382.2
$ws.Range("M31").Value = -5052531.5
$ws.Range("N31").Value = -2972.2

$ws.Range("H34").Value = 3474562.5
$ws.Range("I34").Value = 5052826.5
$ws.Range("J34").Value = 2382.2
$ws.Range("K34").Value = 5052826.5
$ws.Range("L34").Value = 2382.2
$ws.Range("M34").Value = -5052624.5
$ws.Range("N34").Value = -2786.2

$ws.Range("H99").Value = 24998.555
$ws.Range("I99").Value = 35831.168
$ws.Range("K99").Value = 35831.168
$ws.Range("M99").Value = -34333.168

$ws.Range("H126").Value = 24998.555
$ws.Range("I126").Value = 35831.168
$ws.Range("K126").Value = 107493.504
$ws.Range("M126").Value = -105023.504

$ws.Range("H132").Value = 1663.9807
$ws.Range("I132").Value = 1432.3478
$ws.Range("K132").Value = 4297.0434
$ws.Range("M132").Value = -1767.0434

$ws.Range("H134").Value = 1250.6666
$ws.Range("I134").Value = 1121.3793
$ws.Range("K134").Value = 3364.1379
$ws.Range("M134").Value = -829.1379000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4857.5
$ws.Range("I32").Value = 4750
$ws.Range("J32").Value = 4893.3335
$ws.Range("K32").Value = 14250
$ws.Range("L32").Value = 14680.0005
$ws.Range("M32").Value = -13967
$ws.Range("N32").Value = -15246.0005

$ws.Range("H44").Value = 1451.8
$ws.Range("I44").Value = 89.75
$ws.Range("J44").Value = 6900
$ws.Range("K44").Value = 269.25
$ws.Range("L44").Value = 20700
$ws.Range("M44").Value = 128.75
$ws.Range("N44").Value = -21496

$ws.Range("H50").Value = 148330.56
$ws.Range("I50").Value = 1169.9615
$ws.Range("J50").Value = 626602.5
$ws.Range("K50").Value = 3509.8845
$ws.Range("L50").Value = 1879807.5
$ws.Range("M50").Value = -3028.8845
$ws.Range("N50").Value = -1880769.5

$ws.Range("H53").Value = 148330.56
$ws.Range("I53").Value = 1169.9615
$ws.Range("J53").Value = 626602.5
$ws.Range("K53").Value = 3509.8845
$ws.Range("L53").Value = 1879807.5
$ws.Range("M53").Value = -3028.8845
$ws.Range("N53").Value = -1880769.5

$ws.Range("H63").Value = 2037
$ws.Range("I63").Value = 2037
$ws.Range("K63").Value = 6111
$ws.Range("M63").Value = -5362

$ws.Range("H64").Value = 7333
$ws.Range("I64").Value = 6999.5
$ws.Range("K64").Value = 20998.5
$ws.Range("M64").Value = -20728.5

$ws.Range("H66").Value = 2037
$ws.Range("I66").Value = 2037
$ws.Range("K66").Value = 18333
$ws.Range("M66").Value = -14589

$ws.Range("H67").Value = 7333
$ws.Range("I67").Value = 6999.5
$ws.Range("K67").Value = 20998.5
$ws.Range("M67").Value = -20062.5

$ws.Range("H97").Value = 543.2222
$ws.Range("I97").Value = 171.5
$ws.Range("J97").Value = 649.4286
$ws.Range("K97").Value = 514.5
$ws.Range("L97").Value = 1948.2858
$ws.Range("M97").Value = -18.5
$ws.Range("N97").Value = -2940.2858

$ws.Range("H107").Value = 686.6923
$ws.Range("J107").Value = 804.1667
$ws.Range("L107").Value = 2412.5001
$ws.Range("N107").Value = -6252.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4023.5
$ws.Range("I46").Value = 1431
$ws.Range("J46").Value = 4455.5835
$ws.Range("K46").Value = 1431
$ws.Range("L46").Value = 4455.5835
$ws.Range("M46").Value = -1243
$ws.Range("N46").Value = -4831.5835

$ws.Range("H68").Value = 1700
$ws.Range("I68").Value = 1700
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1700
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -951
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1700
$ws.Range("I71").Value = 1700
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8500
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -4756
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 1292.1666
$ws.Range("I82").Value = 1592
$ws.Range("J82").Value = 692.5
$ws.Range("K82").Value = 1592
$ws.Range("L82").Value = 692.5
$ws.Range("M82").Value = -1231
$ws.Range("N82").Value = -1414.5

$ws.Range("H85").Value = 1292.1666
$ws.Range("I85").Value = 1592
$ws.Range("J85").Value = 692.5
$ws.Range("K85").Value = 1592
$ws.Range("L85").Value = 692.5
$ws.Range("M85").Value = -344
$ws.Range("N85").Value = -3188.5

$ws.Range("H132").Value = 1906563.9
$ws.Range("I132").Value = 2779238.5
$ws.Range("J132").Value = 2546.6365
$ws.Range("K132").Value = 8337715.5
$ws.Range("L132").Value = 7639.9095
$ws.Range("M132").Value = -8335185.5
$ws.Range("N132").Value = -12699.9095

$ws.Range("H136").Value = 8338173.5
$ws.Range("I136").Value = 5001359.5
$ws.Range("J136").Value = 25022242
$ws.Range("K136").Value = 15004078.5
$ws.Range("L136").Value = 75066726
$ws.Range("M136").Value = -15001528.5
$ws.Range("N136").Value = -75071826

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2539.6
$ws.Range("I81").Value = 2539.6
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5079.2
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4018.2
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 2539.6
$ws.Range("I84").Value = 2539.6
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 25396
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -20092
$ws.Range("N84").ClearContents()
